# Automatische test-sync: 2025-08-19 19:56:50
# Append a new log row to the "Logs" sheet and bump the matching count on
# the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row of data appended at the bottom of the Logs sheet.
$newRow = 14
$logs.Cells.Item($newRow, 1).Value = "Opvolging retour"
$logs.Cells.Item($newRow, 2).Value = "kwaliteit@testbedrijf123.nl"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 6).Value = "2025-08-19 19:56:38"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Update the summary count on the Dashboard sheet.
$dashboard.Range("B2").Value = 13

# The conditional formatting applied to columns D, G, H, I and J covers the
# data rows of the log; extend those ranges to include the newly added row.
$columns = "D", "G", "H", "I", "J"
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col`2:$col`13")
    $newRange = $logs.Range("$col`2:$col`14")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}
